$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "04/08/2021"
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 20000
$ws.Range("L2").Value = 22000
$ws.Range("M2").Value = 21000
$ws.Range("O2").Value = "Región del Maule"
$ws.Range("P2").Value = 840

$ws.Range("D3").Value = "12/29/2020"
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 30000
$ws.Range("L3").Value = 32000
$ws.Range("M3").Value = 31000
$ws.Range("O3").Value = "Región Metropolitana"
$ws.Range("P3").Value = 1240

$ws.Range("D4").Value = "02/09/2021"
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 25000
$ws.Range("L4").Value = 26000
$ws.Range("M4").Value = 25500
$ws.Range("O4").Value = "Región del Maule"
$ws.Range("P4").Value = 1020

$ws.Range("D5").Value = "06/02/2021"
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 30000
$ws.Range("L5").Value = 32000
$ws.Range("M5").Value = 31000
$ws.Range("O5").Value = "Región Metropolitana"
$ws.Range("P5").Value = 1240

$ws.Range("D6").Value = "04/27/2021"
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 30000
$ws.Range("L6").Value = 32000
$ws.Range("M6").Value = 31000
$ws.Range("O6").Value = "Región Metropolitana"
$ws.Range("P6").Value = 1240

$ws.Range("D7").Value = "12/23/2020"
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 42000
$ws.Range("L7").Value = 44000
$ws.Range("M7").Value = 43000
$ws.Range("O7").Value = "Región de O'Higgins"
$ws.Range("P7").Value = 1720

$ws.Range("D8").Value = "02/03/2021"
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 35000
$ws.Range("L8").Value = 36000
$ws.Range("M8").Value = 35500
$ws.Range("O8").Value = "Región del Maule"
$ws.Range("P8").Value = 1420

$ws.Range("D9").Value = "03/17/2021"
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 22000
$ws.Range("L9").Value = 24000
$ws.Range("M9").Value = 23000
$ws.Range("O9").Value = "Región del Maule"
$ws.Range("P9").Value = 920

$ws.Range("D11").Value = "01/20/2021"
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 26000
$ws.Range("L11").Value = 28000
$ws.Range("M11").Value = 27000
$ws.Range("O11").Value = "Región del Maule"
$ws.Range("P11").Value = 1080

$ws.Range("D12").Value = "05/26/2021"
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 28000
$ws.Range("L12").Value = 30000
$ws.Range("M12").Value = 29000
$ws.Range("O12").Value = "Región Metropolitana"
$ws.Range("P12").Value = 1160

$ws.Range("D13").Value = "05/12/2021"
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 32000
$ws.Range("L13").Value = 34000
$ws.Range("M13").Value = 33000
$ws.Range("O13").Value = "Región Metropolitana"
$ws.Range("P13").Value = 1320

$ws.Range("D14").Value = "01/14/2021"
$ws.Range("J14").Value = 100
$ws.Range("K14").Value = 32000
$ws.Range("L14").Value = 34000
$ws.Range("M14").Value = 33000
$ws.Range("O14").Value = "Región del Maule"
$ws.Range("P14").Value = 1320

$ws.Range("D15").Value = "03/10/2021"
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 22000
$ws.Range("L15").Value = 24000
$ws.Range("M15").Value = 23000
$ws.Range("O15").Value = "Región del Maule"
$ws.Range("P15").Value = 920

$ws.Range("D16").Value = "04/22/2021"
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 28000
$ws.Range("L16").Value = 30000
$ws.Range("M16").Value = 29000
$ws.Range("O16").Value = "Región del Maule"
$ws.Range("P16").Value = 1160

$ws.Range("D17").Value = "02/24/2021"
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 27000
$ws.Range("L17").Value = 28000
$ws.Range("M17").Value = 27500
$ws.Range("O17").Value = "Región Metropolitana"
$ws.Range("P17").Value = 1100

$ws.Range("D18").Value = "04/30/2021"
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 26000
$ws.Range("L18").Value = 27000
$ws.Range("M18").Value = 26500
$ws.Range("O18").Value = "Región Metropolitana"
$ws.Range("P18").Value = 1060

$ws.Range("D19").Value = "03/24/2021"
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 28000
$ws.Range("L19").Value = 30000
$ws.Range("M19").Value = 29000
$ws.Range("O19").Value = "Región del Maule"
$ws.Range("P19").Value = 1160

$ws.Range("D20").Value = "01/07/2021"
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = 25000
$ws.Range("L20").Value = 26000
$ws.Range("M20").Value = 25500
$ws.Range("O20").Value = "Región de O'Higgins"
$ws.Range("P20").Value = 1020

$ws.Range("D21").Value = "02/17/2021"
$ws.Range("J21").Value = 100
$ws.Range("K21").Value = 25000
$ws.Range("L21").Value = 26000
$ws.Range("M21").Value = 25500
$ws.Range("O21").Value = "Región del Maule"
$ws.Range("P21").Value = 1020

